$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 19.92674333333333
$ws.Range("N2").Value = 59.78023
$ws.Range("O2").Value = 0.3447897148135736
$ws.Range("P2").Value = 0.3447897148135735
$ws.Range("Q2").Value = 36.01957460708556
$ws.Range("R2").Value = 324.1761714637701
$ws.Range("S2").Value = 0.3447897148135736
$ws.Range("T2").Value = 0.3447897148135735

# Row 3
$ws.Range("O3").Value = 0.5793221821339875
$ws.Range("P3").Value = 0.5793221821339873
$ws.Range("S3").Value = 0.5793221821339875
$ws.Range("T3").Value = 0.5793221821339873

# Row 4
$ws.Range("M4").Value = 4.385869666666667
$ws.Range("N4").Value = 13.157609
$ws.Range("O4").Value = 0.07588810305243907
$ws.Range("P4").Value = 0.07588810305243905
$ws.Range("Q4").Value = 7.927896547510112
$ws.Range("R4").Value = 71.35106892759102
$ws.Range("S4").Value = 0.07588810305243907
$ws.Range("T4").Value = 0.07588810305243905
